# Test Cases - Login.xlsx update
# - Fill in "Actual Outcome" column (F) with the real observed result for every
#   test case row (it previously just held a placeholder "-").
# - Flip the "Fail/Pass" verdict (column G) for the three backend-login test
#   cases (rows 10-12) from "Fail" to "Pass", since the backend login tests now
#   succeed.
# - Refresh the view state (zoom level and current selection) to match where
#   the author left off reviewing the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Actual Outcome column now matches the Expected Outcome for every test case.
$ws.Range("F2:F12").Value = "Same as expected outcome."

# Backend login test cases (rows 10-12) now pass.
$ws.Range("G10").Value = "Pass"
$ws.Range("G11").Value = "Pass"
$ws.Range("G12").Value = "Pass"

# Update window/view state: zoomed out a bit, selection left on the
# "Actual Outcome" column for the test-case rows.
$win = $excel.ActiveWindow
$win.Zoom = 55
$ws.Range("F2:F12").Select() | Out-Null
